# Apply "Added some unit changes" commit to the workbook.
#
# Sheet1 = FacilityInfo, Sheet2 = CO2LocationInfo

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("FacilityInfo")
$ws2 = $wb.Worksheets.Item("CO2LocationInfo")

# ---------------------------------------------------------------
# FacilityInfo (sheet1) header row relabeling:
# The tkm-* usage columns are reshuffled - their meaning rotates.
# AF: tkm-SZMUsage -> tkm-N1Usage
# AG: pkmUsage     -> tkm-N2Usage
# AH: tkm-N1Usage  -> tkm-N3Usage
# AI: tkm-N2Usage  -> tkm-SZMUsage
# AJ: tkm-N3Usage  -> pkmUsage
# ---------------------------------------------------------------
$ws1.Range("AF1").Value = "tkm-N1Usage"
$ws1.Range("AG1").Value = "tkm-N2Usage"
$ws1.Range("AH1").Value = "tkm-N3Usage"
$ws1.Range("AI1").Value = "tkm-SZMUsage"
$ws1.Range("AJ1").Value = "pkmUsage"

# ---------------------------------------------------------------
# FacilityInfo (sheet1) data row 2 - unit-converted values
# ---------------------------------------------------------------
$ws1.Range("B2").Value  = 4650550.598448258    # Total Cost
$ws1.Range("F2").Value  = 343.481450724496     # hydrogen

$ws1.Range("N2").Value  = 288.5244186085766    # MtDProduction
$ws1.Range("O2").Value  = 288.5244186085766    # MtD-diesel
$ws1.Range("P2").Value  = 0                    # MtGProduction
$ws1.Range("Q2").Value  = 0                    # MtG-gasoline
$ws1.Range("R2").Value  = 288.4396604831751    # PVGasProduction
$ws1.Range("S2").Value  = 288.4396604831752    # PVGas-pkm
$ws1.Range("T2").Value  = 561.5603395168249    # PVDieselProduction
$ws1.Range("U2").Value  = 561.5603395168248    # PVDiesel-pkm

$ws1.Range("AD2").Value = 365.3569032786885    # GasolineHubUsage
$ws1.Range("AE2").Value = 1246.309646477429    # DieselHubUsage

# Values follow their relabeled headers (row values carried along with
# the column's new meaning):
$ws1.Range("AF2").Value = 7.5      # now tkm-N1Usage
$ws1.Range("AG2").Value = 24.2     # now tkm-N2Usage
$ws1.Range("AH2").Value = 130.3    # now tkm-N3Usage
$ws1.Range("AI2").Value = 414.5    # now tkm-SZMUsage
$ws1.Range("AJ2").Value = 850      # now pkmUsage (unit changed from 850008)

# ---------------------------------------------------------------
# CO2LocationInfo (sheet2) data row 2
# ---------------------------------------------------------------
$ws2.Range("B2").Value = "Wacker Chemie AG"    # Name
$ws2.Range("C2").Value = 1612                  # Postal Code
$ws2.Range("D2").Value = 343.481450724496      # Amount Used
